$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.389.09"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.655.39"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.22"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.886.99"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.689.39"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.01"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "27.362.83"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.03"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.46"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.87"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "1.261.51"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.844"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "1.796.93"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.96"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +23.48%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0977"
$ws.Range("E51").Value = "  -0.57%  "
